# Updates the "cryptos" price/volume snapshot table on Sheet1 (columns
# B:Coin, C:Link, D:Price, E:Volume(1h)) to the values captured by the
# latest GitHub Actions run, and re-sorts the NEARProtocol / RenderToken /
# FirstDigitalUSD trio (rows 43-45) into their new rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price cell in this sheet is stored as literal text (e.g. "42.636.25"
# or "1.00"), not a number - Excel would otherwise silently reinterpret
# "1.00" as the number 1. Forcing the cell to Text format before assigning
# the value keeps it as text, and resetting the style back to "Normal"
# afterwards avoids leaving a visible style change on the cell.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '42.636.25'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '2.519.08'
$ws.Range("E3").Value = '  +0.82%  '

Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.10%  '

Set-TextValue "D5" '315.51'
$ws.Range("E5").Value = '  +3.02%  '

Set-TextValue "D6" '95.61'
$ws.Range("E6").Value = '  -0.41%  '

Set-TextValue "D7" '0.575'
$ws.Range("E7").Value = '  -1.60%  '

$ws.Range("E8").Value = '  -0.08%  '

Set-TextValue "D9" '0.533'
$ws.Range("E9").Value = '  -0.35%  '

Set-TextValue "D10" '35.82'
$ws.Range("E10").Value = '  -1.96%  '

Set-TextValue "D11" '0.0810'
$ws.Range("E11").Value = '  -0.10%  '

Set-TextValue "D12" '7.54'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("E13").Value = '  -3.59%  '

$ws.Range("D14").Value = '2.907.56'
$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("D15").Value = '2.502.43'
$ws.Range("E15").Value = '  -0.35%  '

Set-TextValue "D16" '15.25'
$ws.Range("E16").Value = '  -0.58%  '

Set-TextValue "D17" '0.852'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Value = '42.754.21'
$ws.Range("E18").Value = '  +1.00%  '

Set-TextValue "D19" '12.85'
$ws.Range("E19").Value = '  -0.83%  '

Set-TextValue "D20" '6.71'
$ws.Range("E20").Value = '  +4.49%  '

$ws.Range("D21").Value = '0.0₃0959'
$ws.Range("E21").Value = '  -1.30%  '

Set-TextValue "D22" '69.62'
$ws.Range("E22").Value = '  -2.36%  '

Set-TextValue "D23" '249.77'
$ws.Range("E23").Value = '  -1.35%  '

Set-TextValue "D24" '2.95'
$ws.Range("E24").Value = '  +1.14%  '

$ws.Range("E25").Value = '  +2.52%  '

Set-TextValue "D26" '26.58'
$ws.Range("E26").Value = '  -1.21%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("E28").Value = '  +4.00%  '

Set-TextValue "D29" '41.52'
$ws.Range("E29").Value = '  +11.56%  '

Set-TextValue "D30" '10.29'
$ws.Range("E30").Value = '  +1.72%  '

Set-TextValue "D31" '5.95'
$ws.Range("E31").Value = '  +0.17%  '

Set-TextValue "D32" '157.65'
$ws.Range("E32").Value = '  +2.20%  '

Set-TextValue "D33" '2.14'
$ws.Range("E33").Value = '  +3.81%  '

$ws.Range("E34").Value = '  +0.93%  '

Set-TextValue "D35" '2.69'
$ws.Range("E35").Value = '  +3.30%  '

Set-TextValue "D36" '3.30'
$ws.Range("E36").Value = '  +1.17%  '

Set-TextValue "D37" '0.0779'
$ws.Range("E37").Value = '  -0.63%  '

$ws.Range("E38").Value = '  -2.19%  '

Set-TextValue "D39" '0.118'
$ws.Range("E39").Value = '  -0.69%  '

Set-TextValue "D40" '23.38'
$ws.Range("E40").Value = '  -2.06%  '

$ws.Range("E41").Value = '  +16.71%  '

$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D43" '1.00'
$ws.Range("E43").Value = '  +0.38%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D44" '3.33'
$ws.Range("E44").Value = '  -1.62%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D45" '3.79'
$ws.Range("E45").Value = '  -2.17%  '

$ws.Range("D46").Value = '2.038.55'
$ws.Range("E46").Value = '  +0.39%  '

Set-TextValue "D47" '84.31'
$ws.Range("E47").Value = '  -0.10%  '

Set-TextValue "D48" '8.91'
$ws.Range("E48").Value = '  -0.95%  '

Set-TextValue "D49" '75.29'
$ws.Range("E49").Value = '  +3.42%  '

Set-TextValue "D50" '105.52'
$ws.Range("E50").Value = '  +4.53%  '

$ws.Range("D51").Value = '2.764.18'
$ws.Range("E51").Value = '  +0.62%  '
